$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 9366.789000000001
$ws.Range("I28").Value = 1062.6428
$ws.Range("J28").Value = 32618.4
$ws.Range("K28").Value = 1062.6428
$ws.Range("L28").Value = 32618.4
$ws.Range("M28").Value = -577.6428000000001
$ws.Range("N28").Value = -33588.4

$ws.Range("H115").Value = 3322.0833
$ws.Range("I115").Value = 2837.8572
$ws.Range("J115").Value = 4000
$ws.Range("K115").Value = 8513.571599999999
$ws.Range("L115").Value = 12000
$ws.Range("M115").Value = -6946.571599999999
$ws.Range("N115").Value = -15134

$ws.Range("H131").Value = 1558.6154
$ws.Range("I131").Value = 961.6667
$ws.Range("J131").Value = 2070.2856
$ws.Range("K131").Value = 2885.0001
$ws.Range("L131").Value = 6210.8568
$ws.Range("M131").Value = 2154.9999
$ws.Range("N131").Value = -16290.8568

$ws.Range("H137").Value = 2382877
$ws.Range("I137").Value = 3789282.5
$ws.Range("J137").Value = 2806.3845
$ws.Range("K137").Value = 11367847.5
$ws.Range("L137").Value = 8419.1535
$ws.Range("M137").Value = -11365297.5
$ws.Range("N137").Value = -13519.1535

$ws.Range("H140").Value = 76203.19
$ws.Range("J140").Value = 76203.19
$ws.Range("L140").Value = 76203.19
$ws.Range("N140").Value = -86563.19

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1651.1875
$ws.Range("I2").Value = 1594.2142
$ws.Range("J2").Value = 2050
$ws.Range("K2").Value = 1594.2142
$ws.Range("L2").Value = 2050
$ws.Range("M2").Value = -1481.2142
$ws.Range("N2").Value = -2276

$ws.Range("H74").Value = 9616877
$ws.Range("I74").Value = 874.1111
$ws.Range("K74").Value = 874.1111
$ws.Range("M74").Value = -0.1110999999999649

$ws.Range("H77").Value = 9616877
$ws.Range("I77").Value = 874.1111
$ws.Range("K77").Value = 4370.555499999999
$ws.Range("M77").Value = -2.555499999999483

$ws.Range("H116").Value = 1651.1875
$ws.Range("I116").Value = 1594.2142
$ws.Range("J116").Value = 2050
$ws.Range("K116").Value = 1594.2142
$ws.Range("L116").Value = 2050
$ws.Range("M116").Value = 699.7858000000001
$ws.Range("N116").Value = -6638

$ws.Range("H132").Value = 1133963.9
$ws.Range("I132").Value = 1980.3962
$ws.Range("J132").Value = 5133639
$ws.Range("K132").Value = 5941.188599999999
$ws.Range("L132").Value = 15400917
$ws.Range("M132").Value = -3411.188599999999
$ws.Range("N132").Value = -15405977

$ws.Range("H139").Value = 67000.164
$ws.Range("J139").Value = 67000.164
$ws.Range("L139").Value = 67000.164
$ws.Range("N139").Value = -77280.164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1651.1875
$ws.Range("I3").Value = 1594.2142
$ws.Range("J3").Value = 2050
$ws.Range("K3").Value = 1594.2142
$ws.Range("L3").Value = 2050
$ws.Range("M3").Value = -1480.2142
$ws.Range("N3").Value = -2278

$ws.Range("H75").Value = 20471.824
$ws.Range("I75").Value = 3375
$ws.Range("J75").Value = 25732.385
$ws.Range("K75").Value = 3375
$ws.Range("L75").Value = 25732.385
$ws.Range("M75").Value = -2439
$ws.Range("N75").Value = -27604.385

$ws.Range("H78").Value = 20471.824
$ws.Range("I78").Value = 3375
$ws.Range("J78").Value = 25732.385
$ws.Range("K78").Value = 10125
$ws.Range("L78").Value = 77197.155
$ws.Range("M78").Value = -5445
$ws.Range("N78").Value = -86557.155

$ws.Range("H86").Value = 1993.3684
$ws.Range("I86").Value = 1992.5883
$ws.Range("K86").Value = 1992.5883
$ws.Range("M86").Value = -869.5882999999999

$ws.Range("H89").Value = 1993.3684
$ws.Range("I89").Value = 1992.5883
$ws.Range("K89").Value = 9962.941499999999
$ws.Range("M89").Value = -4346.941499999999

$ws.Range("H94").Value = 471.21054
$ws.Range("I94").Value = 466.4375
$ws.Range("J94").Value = 496.66666
$ws.Range("K94").Value = 466.4375
$ws.Range("L94").Value = 496.66666
$ws.Range("M94").Value = -15.4375
$ws.Range("N94").Value = -1398.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7899.028
$ws.Range("I31").Value = 3675.5557
$ws.Range("J31").Value = 9306.852000000001
$ws.Range("K31").Value = 3675.5557
$ws.Range("L31").Value = 9306.852000000001
$ws.Range("M31").Value = -3380.5557
$ws.Range("N31").Value = -9896.852000000001

$ws.Range("H34").Value = 7899.028
$ws.Range("I34").Value = 3675.5557
$ws.Range("J34").Value = 9306.852000000001
$ws.Range("K34").Value = 3675.5557
$ws.Range("L34").Value = 9306.852000000001
$ws.Range("M34").Value = -3473.5557
$ws.Range("N34").Value = -9710.852000000001

$ws.Range("H86").Value = 3342.25
$ws.Range("I86").Value = 3400.7778
$ws.Range("J86").Value = 3166.6667
$ws.Range("K86").Value = 3400.7778
$ws.Range("L86").Value = 3166.6667
$ws.Range("M86").Value = -2277.7778
$ws.Range("N86").Value = -5412.6667

$ws.Range("H89").Value = 3342.25
$ws.Range("I89").Value = 3400.7778
$ws.Range("J89").Value = 3166.6667
$ws.Range("K89").Value = 17003.889
$ws.Range("L89").Value = 15833.3335
$ws.Range("M89").Value = -11387.889
$ws.Range("N89").Value = -27065.3335

$ws.Range("H99").Value = 2245.5
$ws.Range("I99").Value = 1645.7
$ws.Range("J99").Value = 2518.1365
$ws.Range("K99").Value = 1645.7
$ws.Range("L99").Value = 2518.1365
$ws.Range("M99").Value = -147.7
$ws.Range("N99").Value = -5514.136500000001

$ws.Range("H126").Value = 2245.5
$ws.Range("I126").Value = 1645.7
$ws.Range("J126").Value = 2518.1365
$ws.Range("K126").Value = 4937.1
$ws.Range("L126").Value = 7554.4095
$ws.Range("M126").Value = -2467.1
$ws.Range("N126").Value = -12494.4095

$ws.Range("H134").Value = 5438850
$ws.Range("I134").Value = 6582967.5
$ws.Range("J134").Value = 4292.75
$ws.Range("K134").Value = 19748902.5
$ws.Range("L134").Value = 12878.25
$ws.Range("M134").Value = -19746367.5
$ws.Range("N134").Value = -17948.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5999.8
$ws.Range("I81").Value = 1571.2858
$ws.Range("J81").Value = 16333
$ws.Range("K81").Value = 4713.857400000001
$ws.Range("L81").Value = 48999
$ws.Range("M81").Value = -3590.857400000001
$ws.Range("N81").Value = -51245

$ws.Range("H84").Value = 5999.8
$ws.Range("I84").Value = 1571.2858
$ws.Range("J84").Value = 16333
$ws.Range("K84").Value = 14141.5722
$ws.Range("L84").Value = 146997
$ws.Range("M84").Value = -8525.572200000001
$ws.Range("N84").Value = -158229

$ws.Range("H113").Value = 668.51166
$ws.Range("I113").Value = 661.8095
$ws.Range("J113").Value = 674.9091
$ws.Range("K113").Value = 1985.4285
$ws.Range("L113").Value = 2024.7273
$ws.Range("M113").Value = 184.5715
$ws.Range("N113").Value = -6364.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H92").Value = 2750
$ws.Range("J92").Value = 2750
$ws.Range("L92").Value = 2750
$ws.Range("N92").Value = -6494

$ws.Range("H95").Value = 98344
$ws.Range("J95").Value = 98344
$ws.Range("L95").Value = 98344
$ws.Range("N95").Value = -103836

$ws.Range("H98").Value = 98641.5
$ws.Range("J98").Value = 98641.5
$ws.Range("L98").Value = 98641.5
$ws.Range("N98").Value = -104631.5

$ws.Range("H100").Value = 98355
$ws.Range("J100").Value = 98355
$ws.Range("L100").Value = 98355
$ws.Range("N100").Value = -100519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 200000
$ws.Range("J137").Value = 200000
$ws.Range("L137").Value = 200000
$ws.Range("N137").Value = -210200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9921.5
$ws.Range("J54").Value = 9921.5
$ws.Range("L54").Value = 9921.5
$ws.Range("N54").Value = -10961.5

$ws.Range("H76").Value = 44585.75
$ws.Range("J76").Value = 44585.75
$ws.Range("L76").Value = 44585.75
$ws.Range("N76").Value = -45215.75

$ws.Range("H79").Value = 44585.75
$ws.Range("J79").Value = 44585.75
$ws.Range("L79").Value = 44585.75
$ws.Range("N79").Value = -46769.75

$ws.Range("H81").Value = 4023.72
$ws.Range("I81").Value = 4530.0713
$ws.Range("J81").Value = 3379.2727
$ws.Range("K81").Value = 9060.142599999999
$ws.Range("L81").Value = 6758.5454
$ws.Range("M81").Value = -7999.142599999999
$ws.Range("N81").Value = -8880.545399999999

$ws.Range("H84").Value = 4023.72
$ws.Range("I84").Value = 4530.0713
$ws.Range("J84").Value = 3379.2727
$ws.Range("K84").Value = 45300.713
$ws.Range("L84").Value = 33792.727
$ws.Range("M84").Value = -39996.713
$ws.Range("N84").Value = -44400.727
